$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at the top (row 1), shifting all existing rows down by one.
$ws.Rows.Item(1).Insert()

# Set the selection to match the target state.
$ws.Range("A6").Select()
